# "new new 21 maggio" -- add three new withdrawal rows (25-27) to Sheet1,
# matching the existing table layout/formatting, and the four new
# Materiale/UnitaMisura shared strings they introduce.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param([int]$Row, [double]$DateSerial, [string]$Prelevante, [string]$Materiale, [string]$UnitaMisura, [double]$Quantita)

    $colA = $ws.Cells.Item($Row, 1)
    $colA.HorizontalAlignment = -4142   # xlGeneral
    $colA.VerticalAlignment = -4160     # xlTop
    $colA.NumberFormat = "dd/mm/yyyy"
    $colA.Value = $DateSerial

    $colB = $ws.Cells.Item($Row, 2)
    $colB.HorizontalAlignment = -4142
    $colB.VerticalAlignment = -4160
    $colB.NumberFormat = "@"
    $colB.Value = $Prelevante

    $colC = $ws.Cells.Item($Row, 3)
    $colC.HorizontalAlignment = -4142
    $colC.VerticalAlignment = -4160
    $colC.NumberFormat = "@"
    $colC.Value = $Materiale

    $colD = $ws.Cells.Item($Row, 4)
    $colD.HorizontalAlignment = -4142
    $colD.VerticalAlignment = -4160
    $colD.NumberFormat = "@"
    $colD.Value = $UnitaMisura

    $colE = $ws.Cells.Item($Row, 5)
    $colE.HorizontalAlignment = -4142
    $colE.VerticalAlignment = -4160
    $colE.Value = $Quantita
}

Set-Row 25 43239 "Segreteria" "Buste porta documenti plastica" "N°." 100
Set-Row 26 43239 "Segreteria" "Post-It Piccoli" "blocchetti" 12
Set-Row 27 43241 "Licata Rosa" "Tessuto cotone n.7" "Mt." 3
